$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Edit existing shared-string text ---
# Row 130 col P: "Gradual expansion of testing started" -> "Expansion of testing started"
$ws.Cells.Item(130, 16).Value2 = "Expansion of testing started"

# --- 2. Append 4 new daily rows (135-138), 2020-12-08 .. 2020-12-11 ---
$ws.Cells.Item(135, 1).Value2 = 44173
$ws.Cells.Item(135, 2).Value2 = 4673
$ws.Cells.Item(135, 3).Value2 = 339
$ws.Cells.Item(135, 4).Value2 = 1729
$ws.Cells.Item(135, 5).Value2 = 673
$ws.Cells.Item(135, 6).Value2 = 388
$ws.Cells.Item(135, 7).Value2 = 419
$ws.Cells.Item(135, 8).Value2 = 354
$ws.Cells.Item(135, 9).Value2 = 384
$ws.Cells.Item(135, 10).Value2 = 126
$ws.Cells.Item(135, 11).Value2 = 109
$ws.Cells.Item(135, 12).Value2 = 111
$ws.Cells.Item(135, 13).Value2 = 41
$ws.Cells.Item(135, 15).Formula = "=B135-SUM(C135:N135)"
$ws.Cells.Item(135, 16).Value2 = "Further expansion of testing"

$ws.Cells.Item(136, 1).Value2 = 44174
$ws.Cells.Item(136, 2).Value2 = 4734
$ws.Cells.Item(136, 3).Value2 = 489
$ws.Cells.Item(136, 4).Value2 = 1637
$ws.Cells.Item(136, 5).Value2 = 823
$ws.Cells.Item(136, 6).Value2 = 377
$ws.Cells.Item(136, 7).Value2 = 372
$ws.Cells.Item(136, 8).Value2 = 381
$ws.Cells.Item(136, 9).Value2 = 328
$ws.Cells.Item(136, 10).Value2 = 93
$ws.Cells.Item(136, 11).Value2 = 143
$ws.Cells.Item(136, 12).Value2 = 73
$ws.Cells.Item(136, 13).Value2 = 18
$ws.Cells.Item(136, 15).Formula = "=B136-SUM(C136:N136)"

$ws.Cells.Item(137, 1).Value2 = 44175
$ws.Cells.Item(137, 2).Value2 = 4570
$ws.Cells.Item(137, 3).Value2 = 331
$ws.Cells.Item(137, 4).Value2 = 2067
$ws.Cells.Item(137, 5).Value2 = 708
$ws.Cells.Item(137, 6).Value2 = 325
$ws.Cells.Item(137, 7).Value2 = 352
$ws.Cells.Item(137, 8).Value2 = 214
$ws.Cells.Item(137, 9).Value2 = 272
$ws.Cells.Item(137, 10).Value2 = 72
$ws.Cells.Item(137, 11).Value2 = 133
$ws.Cells.Item(137, 12).Value2 = 75
$ws.Cells.Item(137, 13).Value2 = 21
$ws.Cells.Item(137, 15).Formula = "=B137-SUM(C137:N137)"

$ws.Cells.Item(138, 1).Value2 = 44176
$ws.Cells.Item(138, 2).Value2 = 4146
$ws.Cells.Item(138, 3).Value2 = 313
$ws.Cells.Item(138, 4).Value2 = 1708
$ws.Cells.Item(138, 5).Value2 = 636
$ws.Cells.Item(138, 6).Value2 = 397
$ws.Cells.Item(138, 7).Value2 = 281
$ws.Cells.Item(138, 8).Value2 = 268
$ws.Cells.Item(138, 9).Value2 = 252
$ws.Cells.Item(138, 10).Value2 = 82
$ws.Cells.Item(138, 11).Value2 = 116
$ws.Cells.Item(138, 12).Value2 = 69
$ws.Cells.Item(138, 13).Value2 = 24
$ws.Cells.Item(138, 15).Formula = "=B138-SUM(C138:N138)"

# --- 3. New date cells (column A) use the same YYYY-MM-DD date format as the rows above ---
$ws.Range("A135:A138").NumberFormat = "YYYY\-MM\-DD"

# --- 4. Update the visible window / selection to match where the user ended up editing ---
$ws.Application.ActiveWindow.ScrollRow = 118
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("P138").Select()
